$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 242, shifting existing rows 242:301 down to 243:302.
$ws.Rows.Item(242).Insert()

# Populate the newly inserted row 242 with the new observation. All columns
# mirror the (pre-shift) row 242 record except D (date) and J (price), which
# carry the new values.
$ws.Range("A242").Value = 9
$ws.Range("B242").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C242").Value = "Metropolitana"
$ws.Range("D242").Value = 44782
$ws.Range("E242").Value = 13
$ws.Range("F242").Value = 300000001
$ws.Range("G242").Value = "Rabanito"
$ws.Range("H242").Value = "Sin especificar"
$ws.Range("I242").Value = "Primera"
$ws.Range("J242").Value = 7700
$ws.Range("K242").Value = 2500
$ws.Range("L242").Value = 3000
$ws.Range("M242").Value = 2750
$ws.Range("N242").Value = "$/cien unidades (volumen en unidades)"
$ws.Range("O242").Value = "Provincia de Chacabuco"
$ws.Range("P242").Value = 28
$ws.Range("Q242").Value = 100
$ws.Range("R242").Value = "Hortaliza"

# Match the date-formatted number format used by the rest of column D.
$ws.Range("D242").NumberFormat = $ws.Range("D243").NumberFormat
